# Auto-generated Excel COM-interop script to apply numeric value updates
# to the Halicarnassus_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 4
$ws.Range("H4").Value = 1587.6
$ws.Range("I4").Value = 1720.7778
$ws.Range("K4").Value = 1720.7778
$ws.Range("M4").Value = -1606.7778
# ALC row 15
$ws.Range("H15").Value = 2804.6086
$ws.Range("I15").Value = 2804.6086
$ws.Range("K15").Value = 8413.825800000001
$ws.Range("M15").Value = -8244.825800000001
# ALC row 55
$ws.Range("H55").Value = 44.833332
$ws.Range("I55").Value = 51.333332
$ws.Range("K55").Value = 51.333332
$ws.Range("M55").Value = 162.666668
# ALC row 86
$ws.Range("H86").Value = 3288.8572
$ws.Range("J86").Value = 3462.0833
$ws.Range("L86").Value = 3462.0833
$ws.Range("N86").Value = -5708.0833
# ALC row 89
$ws.Range("H89").Value = 3288.8572
$ws.Range("J89").Value = 3462.0833
$ws.Range("L89").Value = 17310.4165
$ws.Range("N89").Value = -28542.4165
# ALC row 107
$ws.Range("H107").Value = 234.42308
$ws.Range("I107").Value = 97.95238000000001
$ws.Range("K107").Value = 97.95238000000001
$ws.Range("M107").Value = 1822.04762
# ALC row 121
$ws.Range("H121").Value = 903.2
$ws.Range("J121").Value = 903.2
$ws.Range("L121").Value = 2709.6
$ws.Range("N121").Value = -6203.6
# ALC row 132
$ws.Range("H132").Value = 2463
$ws.Range("I132").Value = 1336.2
$ws.Range("K132").Value = 4008.6
$ws.Range("M132").Value = -1478.6
# ALC row 137
$ws.Range("H137").Value = 4664.143
$ws.Range("I137").Value = 2139.111
$ws.Range("K137").Value = 6417.333
$ws.Range("M137").Value = -3867.333
# ALC row 140
$ws.Range("H140").Value = 90780
$ws.Range("J140").Value = 90780
$ws.Range("L140").Value = 90780
$ws.Range("N140").Value = -101140
# ALC row 141
$ws.Range("H141").Value = 2090.9167
$ws.Range("I141").Value = 1199
$ws.Range("K141").Value = 3597
$ws.Range("M141").Value = 1583

$ws = $wb.Worksheets.Item("ARM")
# ARM row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
# ARM row 32
$ws.Range("H32").Value = 12460.267
$ws.Range("I32").Value = 9608.434999999999
$ws.Range("K32").Value = 9608.434999999999
$ws.Range("M32").Value = -9321.434999999999
# ARM row 61
$ws.Range("H61").Value = 4462
$ws.Range("I61").Value = 4462
$ws.Range("K61").Value = 4462
$ws.Range("M61").Value = -4250
# ARM row 74
$ws.Range("H74").Value = 2339.4
$ws.Range("I74").Value = 1862.5
$ws.Range("J74").Value = 4247
$ws.Range("K74").Value = 1862.5
$ws.Range("L74").Value = 4247
$ws.Range("M74").Value = -988.5
$ws.Range("N74").Value = -5995
# ARM row 77
$ws.Range("H77").Value = 2339.4
$ws.Range("I77").Value = 1862.5
$ws.Range("J77").Value = 4247
$ws.Range("K77").Value = 9312.5
$ws.Range("L77").Value = 21235
$ws.Range("M77").Value = -4944.5
$ws.Range("N77").Value = -29971
# ARM row 136
$ws.Range("H136").Value = 4462
$ws.Range("I136").Value = 4462
$ws.Range("K136").Value = 13386
$ws.Range("M136").Value = -10836

$ws = $wb.Worksheets.Item("BSM")
# BSM row 64
$ws.Range("H64").Value = 581.8333
$ws.Range("I64").Value = 314.66666
$ws.Range("K64").Value = 314.66666
$ws.Range("M64").Value = -89.66665999999998
# BSM row 67
$ws.Range("H67").Value = 581.8333
$ws.Range("I67").Value = 314.66666
$ws.Range("K67").Value = 314.66666
$ws.Range("M67").Value = 465.33334
# BSM row 80
$ws.Range("H80").Value = 227.26315
$ws.Range("I80").Value = 173.16667
$ws.Range("J80").Value = 252.23077
$ws.Range("K80").Value = 173.16667
$ws.Range("L80").Value = 252.23077
$ws.Range("M80").Value = 824.8333299999999
$ws.Range("N80").Value = -2248.23077
# BSM row 83
$ws.Range("H83").Value = 227.26315
$ws.Range("I83").Value = 173.16667
$ws.Range("J83").Value = 252.23077
$ws.Range("K83").Value = 865.8333500000001
$ws.Range("L83").Value = 1261.15385
$ws.Range("M83").Value = 4126.16665
$ws.Range("N83").Value = -11245.15385
# BSM row 86
$ws.Range("H86").Value = 4788.4585
$ws.Range("I86").Value = 3389.647
$ws.Range("J86").Value = 8185.5713
$ws.Range("K86").Value = 3389.647
$ws.Range("L86").Value = 8185.5713
$ws.Range("M86").Value = -2266.647
$ws.Range("N86").Value = -10431.5713
# BSM row 89
$ws.Range("H89").Value = 4788.4585
$ws.Range("I89").Value = 3389.647
$ws.Range("J89").Value = 8185.5713
$ws.Range("K89").Value = 16948.235
$ws.Range("L89").Value = 40927.85649999999
$ws.Range("M89").Value = -11332.235
$ws.Range("N89").Value = -52159.85649999999
# BSM row 94
$ws.Range("H94").Value = 1037.2
$ws.Range("I94").Value = 996.5
$ws.Range("K94").Value = 996.5
$ws.Range("M94").Value = -545.5
# BSM row 134
$ws.Range("H134").Value = 1812.8182
$ws.Range("I134").Value = 882.44446
$ws.Range("K134").Value = 2647.33338
$ws.Range("M134").Value = -112.33338

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 5803.3076
$ws.Range("I31").Value = 4831.1
$ws.Range("K31").Value = 4831.1
$ws.Range("M31").Value = -4536.1
# CRP row 34
$ws.Range("H34").Value = 5803.3076
$ws.Range("I34").Value = 4831.1
$ws.Range("K34").Value = 4831.1
$ws.Range("M34").Value = -4629.1
# CRP row 107
$ws.Range("H107").Value = 409.41666
$ws.Range("I107").Value = 329.8421
$ws.Range("K107").Value = 329.8421
$ws.Range("M107").Value = 1590.1579

$ws = $wb.Worksheets.Item("CUL")
# CUL row 122
$ws.Range("H122").Value = 949.3333
$ws.Range("J122").Value = 949.3333
$ws.Range("L122").Value = 8543.9997
$ws.Range("N122").Value = -13443.9997
# CUL row 141
$ws.Range("H141").Value = 1991.3636
$ws.Range("I141").Value = 1890.5
$ws.Range("K141").Value = 5671.5
$ws.Range("M141").Value = -491.5

$ws = $wb.Worksheets.Item("GSM")
# GSM row 55
$ws.Range("H55").Value = 7124.4
$ws.Range("I55").Value = 5905.5
$ws.Range("K55").Value = 5905.5
$ws.Range("M55").Value = -5578.5
# GSM row 70
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 10000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -10540
# GSM row 73
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 10000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -11872

$ws = $wb.Worksheets.Item("LTW")
# LTW row 132
$ws.Range("H132").Value = 7733.4707
$ws.Range("I132").Value = 5748.6
$ws.Range("K132").Value = 17245.8
$ws.Range("M132").Value = -14715.8
# LTW row 136
$ws.Range("H136").Value = 7574.5
$ws.Range("I136").Value = 6766
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 20298
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -17748
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("WVR")
# WVR row 2
$ws.Range("H2").Value = 3056.25
$ws.Range("I2").Value = 3056.25
$ws.Range("K2").Value = 3056.25
$ws.Range("M2").Value = -2944.25
# WVR row 132
$ws.Range("H132").Value = 4561.0293
$ws.Range("I132").Value = 4150.4443
$ws.Range("K132").Value = 12451.3329
$ws.Range("M132").Value = -9921.332900000001
# WVR row 136
$ws.Range("H136").Value = 5145
$ws.Range("I136").Value = 3832.44
$ws.Range("K136").Value = 11497.32
$ws.Range("M136").Value = -8947.32
